$d = $word.ActiveDocument

# 1. Add the new "Abstract Title" paragraph style (styleId "AbstractTitle"),
#    based on Normal, followed by Abstract, centered, bold, small blue-grey text.
$abstractTitle = $d.Styles.Add("Abstract Title", 1)
$abstractTitle.BaseStyle = $d.Styles("Normal")
$abstractTitle.NextParagraphStyle = $d.Styles("Abstract")
$abstractTitle.QuickStyle = $true
$abstractTitle.ParagraphFormat.KeepWithNext = $true
$abstractTitle.ParagraphFormat.KeepTogether = $true
$abstractTitle.ParagraphFormat.Alignment = 1
$abstractTitle.ParagraphFormat.SpaceBefore = 15
$abstractTitle.ParagraphFormat.SpaceAfter = 0
$abstractTitle.Font.Size = 10
$abstractTitle.Font.SizeBi = 10
$abstractTitle.Font.Bold = $true
$abstractTitle.Font.Color = 0x8A5A34

# 2. Abstract style: reduce space-before from 15pt (300) to 5pt (100).
$abstract = $d.Styles("Abstract")
$abstract.ParagraphFormat.SpaceBefore = 5

# 3. ImportTok character style: add green, bold colouring.
$importTok = $d.Styles("ImportTok")
$importTok.Font.Color = 0x008000
$importTok.Font.Bold = $true

# 4. BuiltInTok character style: add green colouring.
$builtInTok = $d.Styles("BuiltInTok")
$builtInTok.Font.Color = 0x008000
